$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "FR79 2282 6763 09PW 6IWO 96PK Y07"
$ws.Range("A3").Value = "FR31 4834 6092 31UZ W7ZT UGFB R22"
$ws.Range("A4").Value = "EE68 7340 1352 5566 3731"
$ws.Range("A5").Value = "GR11 1018 094X IYFC PGT3 5INY 7IC"
$ws.Range("A6").Value = "AZ19 IUST XYBK VDU4 XFZ2 KMY3 1IKE"
$ws.Range("A7").Value = "FR62 0841 1230 28O5 GDVJ MDBC L03"
$ws.Range("A8").Value = "PK73 OMZO C53E 3IYE SXAD RTFN"
$ws.Range("A9").Value = "EE44 8823 9718 1278 6819"
$ws.Range("A10").Value = "CR69 6448 2334 9876 1879 8"
$ws.Range("A11").Value = "NL09 TLMT 4085 0135 17"
$ws.Range("A12").Value = "FR45 5494 2255 49OE YVWD NGW6 Z14"
$ws.Range("A13").Value = "HR20 9159 6225 3934 9298 9"
$ws.Range("A14").Value = "EE49 7319 9421 4974 1039"
$ws.Range("A15").Value = "FR26 5518 8001 92IY JWZX FEJN C44"
$ws.Range("A16").Value = "SA51 71XS OWWW 3KEB N83G HPM9"
$ws.Range("A17").Value = "FI35 1242 1345 9527 18"
$ws.Range("A18").Value = "IT60 Q802 2504 904M WSUP GBFI UA4"
$ws.Range("A19").Value = "HU18 3488 8407 0210 9257 0388 5034"
$ws.Range("A20").Value = "HR68 1204 4287 4917 3008 3"
$ws.Range("A21").Value = "SE80 8624 0816 4622 9352 3647"
$ws.Range("A22").Value = "SI96 2089 1193 9003 267"
$ws.Range("A23").Value = "SE58 0549 2590 2624 9416 9223"
$ws.Range("A24").Value = "FR08 9965 9725 56QB DKF1 0MIN V61"
$ws.Range("A25").Value = "FR44 5480 3794 60TT ZIJN KQQ3 F27"
$ws.Range("A26").Value = "SA02 473T FNEI JBKI OPR0 VG0F"
$ws.Range("A27").Value = "TR93 0116 92B2 WUSO G3LX EHUX PG"
$ws.Range("A28").Value = "SA14 085X WT62 45NY J1MD JU14"
$ws.Range("A29").Value = "FR55 8004 4030 44EB BN7T BDI1 Q19"
$ws.Range("A30").Value = "MR69 1237 3041 8942 8641 5683 215"
$ws.Range("A31").Value = "IE57 ZYQA 5080 9339 3722 31"
$ws.Range("A32").Value = "MT21 GHVM 2935 6HIR ZVLQ 9ZB4 OPJS 5D6"
$ws.Range("A33").Value = "RO21 RAMB 9LQ7 OEPB 5NEI 1MDO"
$ws.Range("A34").Value = "MD47 YNEH Z2QJ UHWR UBFB KIBJ"
$ws.Range("A35").Value = "FR51 7169 3117 31IQ PRHL TDAT L25"
$ws.Range("A36").Value = "FR74 3246 0841 70UY HYFL P81N X77"
$ws.Range("A37").Value = "GR48 8581 0775 PIN0 1XSD AMKQ MUE"
$ws.Range("A38").Value = "AZ86 KOHM JDOH AEB1 AFRG RAYC DTIH"
$ws.Range("A39").Value = "CR49 8336 6285 9676 0875 2"
$ws.Range("A40").Value = "GT71 PN5L L1PP ECLG 8YE1 SART MMEG"
$ws.Range("A41").Value = "AT12 1475 9491 3950 1051"
$ws.Range("A42").Value = "IL83 5273 9999 4201 5858 982"
$ws.Range("A43").Value = "ES05 6582 4291 9678 5877 2373"
$ws.Range("A44").Value = "FR13 1333 7156 37FM INJW 7JM1 E34"
$ws.Range("A45").Value = "RS13 8531 4475 6204 9750 80"
$ws.Range("A46").Value = "BE42 3218 6233 1524"
$ws.Range("A47").Value = "PS21 KPVI N7BL LJS2 H78H AR5J USXB J"
$ws.Range("A48").Value = "MD10 53CZ M45Y GHYV MKKZ ZBMY"
$ws.Range("A49").Value = "EE89 8422 6975 3891 9191"
$ws.Range("A50").Value = "IT78 X612 8332 086X D7XC 8I04 BNS"

$ws.Activate()
$ws.Range("E6").Select()
